$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the data set ("RM 232" at row 26
# and "SC 92" at row 28). Delete the higher-numbered row first so the second
# delete's row index is unaffected.
$ws.Rows("28").Delete()
$ws.Rows("26").Delete()

# After the deletions, the remaining rows have shifted up by two. Apply the
# per-cell value corrections (these are the cells whose imputed/missing
# values changed between the "before" and "after" snapshots).
$ws.Range("D2").Value = -13.5
$ws.Range("D3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("E8").Value = -6.6
$ws.Range("E10").Value = -6.1
$ws.Range("D11").Value = -15.5
$ws.Range("E12").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E15").Value = -8.4
$ws.Range("E18").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("D21").Value = -14.3
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = -7.1
$ws.Range("E27").Value = -10
$ws.Range("C29").Value = ""
$ws.Range("E29").Value = ""
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = ""
